# Scheduled-runner style refresh: update cached market-price / leve-profit
# values across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 4649.5
$ws.Range("J43").Value = 4649.5
$ws.Range("L43").Value = 4649.5
$ws.Range("N43").Value = -4787.5
$ws.Range("H70").Value = 7137.5
$ws.Range("J70").Value = 4516
$ws.Range("L70").Value = 13548
$ws.Range("N70").Value = -14088
$ws.Range("H73").Value = 7137.5
$ws.Range("J73").Value = 4516
$ws.Range("L73").Value = 13548
$ws.Range("N73").Value = -15420
$ws.Range("H80").Value = 1050.091
$ws.Range("I80").Value = 762.75
$ws.Range("J80").Value = 1214.2858
$ws.Range("K80").Value = 2288.25
$ws.Range("L80").Value = 3642.8574
$ws.Range("M80").Value = -1290.25
$ws.Range("N80").Value = -5638.857400000001
$ws.Range("H83").Value = 1050.091
$ws.Range("I83").Value = 762.75
$ws.Range("J83").Value = 1214.2858
$ws.Range("K83").Value = 6864.75
$ws.Range("L83").Value = 10928.5722
$ws.Range("M83").Value = -1872.75
$ws.Range("N83").Value = -20912.5722
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H116").Value = 4592.25
$ws.Range("I116").Value = 4495
$ws.Range("K116").Value = 4495
$ws.Range("M116").Value = -1053

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1602
$ws.Range("I2").Value = 1302.6666
$ws.Range("K2").Value = 1302.6666
$ws.Range("M2").Value = -1189.6666
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H39").Value = 2558
$ws.Range("I39").Value = 2558
$ws.Range("K39").Value = 2558
$ws.Range("M39").Value = -2038
$ws.Range("H41").Value = 1870.6666
$ws.Range("I41").Value = 1870.6666
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1870.6666
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -1456.6666
$ws.Range("H116").Value = 1602
$ws.Range("I116").Value = 1302.6666
$ws.Range("K116").Value = 1302.6666
$ws.Range("M116").Value = 991.3334
$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1499.5
$ws.Range("K122").Value = 4498.5
$ws.Range("M122").Value = -2048.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1602
$ws.Range("I3").Value = 1302.6666
$ws.Range("K3").Value = 1302.6666
$ws.Range("M3").Value = -1188.6666
$ws.Range("H134").Value = 8000
$ws.Range("I134").Value = 8000
$ws.Range("K134").Value = 24000
$ws.Range("M134").Value = -21465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 232.02325
$ws.Range("I7").Value = 265.8
$ws.Range("J7").Value = 221.78787
$ws.Range("K7").Value = 265.8
$ws.Range("L7").Value = 221.78787
$ws.Range("M7").Value = -152.8
$ws.Range("N7").Value = -447.78787
$ws.Range("H25").Value = 340.66666
$ws.Range("I25").Value = 340.66666
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 340.66666
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -166.66666
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H80").Value = 60000
$ws.Range("J80").Value = 60000
$ws.Range("L80").Value = 60000
$ws.Range("N80").Value = -62246
$ws.Range("H83").Value = 60000
$ws.Range("J83").Value = 60000
$ws.Range("L83").Value = 180000
$ws.Range("N83").Value = -191232
$ws.Range("H86").Value = 7623.5
$ws.Range("I86").Value = 5998.3335
$ws.Range("K86").Value = 5998.3335
$ws.Range("M86").Value = -4875.3335
$ws.Range("H88").Value = 18321.5
$ws.Range("J88").Value = 18321.5
$ws.Range("L88").Value = 18321.5
$ws.Range("N88").Value = -19133.5
$ws.Range("H89").Value = 7623.5
$ws.Range("I89").Value = 5998.3335
$ws.Range("K89").Value = 29991.6675
$ws.Range("M89").Value = -24375.6675
$ws.Range("H91").Value = 18321.5
$ws.Range("J91").Value = 18321.5
$ws.Range("L91").Value = 18321.5
$ws.Range("N91").Value = -21129.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").ClearContents()
$ws.Range("N95").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1924.3
$ws.Range("I14").Value = 1924.3
$ws.Range("K14").Value = 5772.9
$ws.Range("M14").Value = -5599.9
$ws.Range("H114").Value = 399.5
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H128").Value = 200000
$ws.Range("I128").Value = 200000
$ws.Range("K128").Value = 600000
$ws.Range("M128").Value = -595020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2800.75
$ws.Range("I3").Value = 601.5
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 601.5
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -485.5
$ws.Range("N3").Value = -5232
$ws.Range("H12").Value = 2000
$ws.Range("J12").Value = 2000
$ws.Range("L12").Value = 2000
$ws.Range("N12").Value = -2280
$ws.Range("H122").Value = 8268.333000000001
$ws.Range("I122").Value = 2403.5
$ws.Range("K122").Value = 7210.5
$ws.Range("M122").Value = -4760.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 10000000
$ws.Range("I24").Value = 10000000
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 10000000
$ws.Range("L24").ClearContents()
$ws.Range("M24").Value = -9999657
$ws.Range("N24").Value = 0
$ws.Range("H122").Value = 3904.6
$ws.Range("J122").Value = 3874.6667
$ws.Range("L122").Value = 11624.0001
$ws.Range("N122").Value = -16524.0001
$ws.Range("H132").Value = 9728.691999999999
$ws.Range("I132").Value = 9773.777
$ws.Range("K132").Value = 29321.331
$ws.Range("M132").Value = -26791.331

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2000.8572
$ws.Range("I113").Value = 1400.75
$ws.Range("K113").Value = 4202.25
$ws.Range("M113").Value = -2032.25
$ws.Range("H126").Value = 3342.7778
$ws.Range("I126").Value = 3512.1428
$ws.Range("K126").Value = 10536.4284
$ws.Range("M126").Value = -8066.428400000001
$ws.Range("H132").Value = 9749.75
$ws.Range("I132").Value = 8000
$ws.Range("J132").Value = 11499.5
$ws.Range("K132").Value = 24000
$ws.Range("L132").Value = 34498.5
$ws.Range("M132").Value = -21470
